$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.680.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.583.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.15%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.47%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.67%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.10%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.808.25"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.561.44"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.63%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.10%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.644.82"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.50%  "

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.64%  "

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.00%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.37%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.48%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.70"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.63%  "

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.76"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.08%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.11"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.91%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.05%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.68%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.63%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.386.54"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.75%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.39%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.59%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.74%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.04%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.976"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.62%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.47"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.97%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.80%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.719.31"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.92"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0974"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0499"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.21%  "
